$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty F10 cell (inlineStr with no content) from the previous edit
$ws.Cells.Item(10, 6).ClearContents()

# Append new account rows (11-27) mirroring the scraped data
$ws.Cells.Item(11, 1).Value = "Sarah Garcia"
$ws.Cells.Item(11, 2).Value = "sarah.brown81@hotmail.com"
$ws.Cells.Item(11, 3).Value = "Pass3292"
$ws.Cells.Item(11, 4).Value = "4 December 1989"
$ws.Cells.Item(11, 5).Value = "male"
$ws.Cells.Item(11, 7).Value = "2025-04-19 20:53:48"

$ws.Cells.Item(12, 1).Value = "Sarah Garcia"
$ws.Cells.Item(12, 2).Value = "sarah.brown81@hotmail.com"
$ws.Cells.Item(12, 3).Value = "Pass3292"
$ws.Cells.Item(12, 4).Value = "4 December 1989"
$ws.Cells.Item(12, 5).Value = "male"
$ws.Cells.Item(12, 7).Value = "2025-04-19 20:53:48"

$ws.Cells.Item(13, 1).Value = "Sarah Jones"
$ws.Cells.Item(13, 2).Value = "sarah.johnson72@outlook.com"
$ws.Cells.Item(13, 3).Value = "Pass3942"
$ws.Cells.Item(13, 4).Value = "28 July 1992"
$ws.Cells.Item(13, 5).Value = "male"
$ws.Cells.Item(13, 7).Value = "2025-04-19 20:54:15"

$ws.Cells.Item(14, 1).Value = "Sarah Jones"
$ws.Cells.Item(14, 2).Value = "sarah.johnson72@outlook.com"
$ws.Cells.Item(14, 3).Value = "Pass3942"
$ws.Cells.Item(14, 4).Value = "28 July 1992"
$ws.Cells.Item(14, 5).Value = "male"
$ws.Cells.Item(14, 7).Value = "2025-04-19 20:54:15"

$ws.Cells.Item(15, 1).Value = "James Brown"
$ws.Cells.Item(15, 2).Value = "james.johnson24@protonmail.com"
$ws.Cells.Item(15, 3).Value = "Pass4590"
$ws.Cells.Item(15, 4).Value = "27 October 1992"
$ws.Cells.Item(15, 5).Value = "female"
$ws.Cells.Item(15, 7).Value = "2025-04-19 20:54:44"

$ws.Cells.Item(16, 1).Value = "James Brown"
$ws.Cells.Item(16, 2).Value = "james.johnson24@protonmail.com"
$ws.Cells.Item(16, 3).Value = "Pass4590"
$ws.Cells.Item(16, 4).Value = "27 October 1992"
$ws.Cells.Item(16, 5).Value = "female"
$ws.Cells.Item(16, 7).Value = "2025-04-19 20:54:48"

$ws.Cells.Item(17, 1).Value = "James Jones"
$ws.Cells.Item(17, 2).Value = "james.williams45@gmail.com"
$ws.Cells.Item(17, 3).Value = "Pass2039"
$ws.Cells.Item(17, 4).Value = "8 September 1985"
$ws.Cells.Item(17, 5).Value = "female"
$ws.Cells.Item(17, 7).Value = "2025-04-19 20:57:00"

$ws.Cells.Item(18, 1).Value = "James Jones"
$ws.Cells.Item(18, 2).Value = "james.williams45@gmail.com"
$ws.Cells.Item(18, 3).Value = "Pass2039"
$ws.Cells.Item(18, 4).Value = "8 September 1985"
$ws.Cells.Item(18, 5).Value = "female"
$ws.Cells.Item(18, 7).Value = "2025-04-19 20:57:00"

$ws.Cells.Item(19, 1).Value = "John Smith"
$ws.Cells.Item(19, 2).Value = "john.johnson6@protonmail.com"
$ws.Cells.Item(19, 3).Value = "Pass2867"
$ws.Cells.Item(19, 4).Value = "22 March 1995"
$ws.Cells.Item(19, 5).Value = "male"
$ws.Cells.Item(19, 7).Value = "2025-04-19 21:07:29"

$ws.Cells.Item(20, 1).Value = "John Smith"
$ws.Cells.Item(20, 2).Value = "john.johnson6@protonmail.com"
$ws.Cells.Item(20, 3).Value = "Pass2867"
$ws.Cells.Item(20, 4).Value = "22 March 1995"
$ws.Cells.Item(20, 5).Value = "male"
$ws.Cells.Item(20, 7).Value = "2025-04-19 21:07:29"

$ws.Cells.Item(21, 1).Value = "Emma Brown"
$ws.Cells.Item(21, 2).Value = "Loading"
$ws.Cells.Item(21, 3).Value = "Pass5452"
$ws.Cells.Item(21, 4).Value = "3 December 1987"
$ws.Cells.Item(21, 5).Value = "male"
$ws.Cells.Item(21, 7).Value = "2025-04-19 21:12:02"

$ws.Cells.Item(22, 1).Value = "Emma Brown"
$ws.Cells.Item(22, 2).Value = "Loading"
$ws.Cells.Item(22, 3).Value = "Pass5452"
$ws.Cells.Item(22, 4).Value = "3 December 1987"
$ws.Cells.Item(22, 5).Value = "male"
$ws.Cells.Item(22, 7).Value = "2025-04-19 21:12:02"

$ws.Cells.Item(23, 1).Value = "John Williams"
$ws.Cells.Item(23, 2).Value = "Loading"
$ws.Cells.Item(23, 3).Value = "Pass8711"
$ws.Cells.Item(23, 4).Value = "4 October 1990"
$ws.Cells.Item(23, 5).Value = "male"
$ws.Cells.Item(23, 7).Value = "2025-04-19 21:12:36"

$ws.Cells.Item(24, 1).Value = "Emma Jones"
$ws.Cells.Item(24, 2).Value = "wasaro6972@agiuse.com"
$ws.Cells.Item(24, 3).Value = "Pass6375"
$ws.Cells.Item(24, 4).Value = "20 September 1990"
$ws.Cells.Item(24, 5).Value = "male"
$ws.Cells.Item(24, 7).Value = "2025-04-19 21:14:16"

$ws.Cells.Item(25, 1).Value = "Emma Jones"
$ws.Cells.Item(25, 2).Value = "wasaro6972@agiuse.com"
$ws.Cells.Item(25, 3).Value = "Pass6375"
$ws.Cells.Item(25, 4).Value = "20 September 1990"
$ws.Cells.Item(25, 5).Value = "male"
$ws.Cells.Item(25, 7).Value = "2025-04-19 21:14:16"

$ws.Cells.Item(26, 1).Value = "John Garcia"
$ws.Cells.Item(26, 2).Value = "nilono4638@f5url.com"
$ws.Cells.Item(26, 3).Value = "Pass3288"
$ws.Cells.Item(26, 4).Value = "20 March 1988"
$ws.Cells.Item(26, 5).Value = "male"
$ws.Cells.Item(26, 7).Value = "2025-04-19 21:15:00"

$ws.Cells.Item(27, 1).Value = "John Garcia"
$ws.Cells.Item(27, 2).Value = "nilono4638@f5url.com"
$ws.Cells.Item(27, 3).Value = "Pass3288"
$ws.Cells.Item(27, 4).Value = "20 March 1988"
$ws.Cells.Item(27, 5).Value = "male"
$ws.Cells.Item(27, 7).Value = "2025-04-19 21:15:00"

Write-Host "Done"
